# Remove the empty paragraph, the "Ver no Jupiter..." paragraph, and the
# "(c) 2020 ..." footer paragraph that follow the
# "LOB1045: Leitura e Produção de Textos Acadêmicos (Requisito)" line,
# while leaving the paragraph that originally trailed the footer
# (and the subsequent page-break paragraph) untouched.

$d = $word.ActiveDocument

# Locate the "LOB1045" requirement paragraph.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "LOB1045*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -ne -1) {
    # The three paragraphs that must disappear are the ones immediately
    # following the anchor: the blank spacer paragraph, the "Ver no
    # Jupiter..." line, and the "(c) 2020..." footer line.
    $firstToRemove = $anchorIndex + 1
    $lastToRemove = $anchorIndex + 3

    $startRange = $d.Paragraphs.Item($firstToRemove).Range.Start
    $endRange = $d.Paragraphs.Item($lastToRemove).Range.End

    $killRange = $d.Range($startRange, $endRange)
    $killRange.Delete()
}
